$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("Q2").Value = "easy:81.62"
$ws.Range("R2").Value = "very:13.31"
$ws.Range("S2").Value = "neutral:4.69"
$ws.Range("U2").Value = "None"

# Row 3 updates
$ws.Range("Q3").Value = "easy:87.43"
$ws.Range("R3").Value = "very:8.17"
$ws.Range("S3").Value = "neutral:3.93"
$ws.Range("U3").Value = "difficult:0.07"
